$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update shifted rows 96-201 (D, J, K, L, M, P columns)
$rowData = @{
    96 = @{ D=44539; J=60; K=6000; L=6000; M=6000; P=167 }
    97 = @{ D=44232; J=180; K=6000; L=6000; M=6000; P=167 }
    98 = @{ D=44411; J=180; K=6000; L=6000; M=6000; P=167 }
    99 = @{ D=44487; J=60; K=5500; L=6000; M=5750; P=160 }
    100 = @{ D=44176; J=180; K=6000; L=6000; M=6000; P=167 }
    101 = @{ D=44165; J=60; K=6000; L=6000; M=6000; P=167 }
    102 = @{ D=44326; J=30; K=6000; L=6000; M=6000; P=167 }
    103 = @{ D=44530; J=180; K=6000; L=6000; M=6000; P=167 }
    104 = @{ D=44336; J=60; K=6000; L=6000; M=6000; P=167 }
    105 = @{ D=44292; J=180; K=6500; L=6500; M=6500; P=181 }
    106 = @{ D=44252; J=90; K=6000; L=6000; M=6000; P=167 }
    107 = @{ D=44400; J=180; K=6000; L=6000; M=6000; P=167 }
    108 = @{ D=44245; J=80; K=6000; L=6000; M=6000; P=167 }
    109 = @{ D=44379; J=180; K=6000; L=6000; M=6000; P=167 }
    110 = @{ D=44243; J=180; K=6000; L=6000; M=6000; P=167 }
    111 = @{ D=44484; J=160; K=5500; L=5500; M=5500; P=153 }
    112 = @{ D=44250; J=180; K=6000; L=6000; M=6000; P=167 }
    113 = @{ D=44363; J=30; K=6000; L=6000; M=6000; P=167 }
    114 = @{ D=44476; J=80; K=6000; L=6000; M=6000; P=167 }
    115 = @{ D=44236; J=180; K=6000; L=6000; M=6000; P=167 }
    116 = @{ D=44354; J=30; K=6000; L=6000; M=6000; P=167 }
    117 = @{ D=44196; J=60; K=6000; L=7000; M=6500; P=181 }
    118 = @{ D=44438; J=70; K=6500; L=6500; M=6500; P=181 }
    119 = @{ D=44306; J=180; K=6000; L=6000; M=6000; P=167 }
    120 = @{ D=44397; J=180; K=6000; L=6000; M=6000; P=167 }
    121 = @{ D=44407; J=180; K=6000; L=6000; M=6000; P=167 }
    122 = @{ D=44285; J=180; K=6000; L=6000; M=6000; P=167 }
    123 = @{ D=44371; J=60; K=6000; L=6000; M=6000; P=167 }
    124 = @{ D=44263; J=30; K=7000; L=7000; M=7000; P=194 }
    125 = @{ D=44390; J=180; K=6000; L=6000; M=6000; P=167 }
    126 = @{ D=44277; J=60; K=6500; L=6500; M=6500; P=181 }
    127 = @{ D=44159; J=150; K=6500; L=6500; M=6500; P=181 }
    128 = @{ D=44355; J=140; K=6000; L=6000; M=6000; P=167 }
    129 = @{ D=44529; J=90; K=6000; L=6000; M=6000; P=167 }
    130 = @{ D=44221; J=80; K=6500; L=6500; M=6500; P=181 }
    131 = @{ D=44417; J=90; K=7000; L=7000; M=7000; P=194 }
    132 = @{ D=44406; J=90; K=6000; L=6000; M=6000; P=167 }
    133 = @{ D=44384; J=27; K=6000; L=6000; M=6000; P=167 }
    134 = @{ D=44432; J=180; K=7000; L=7000; M=7000; P=194 }
    135 = @{ D=44403; J=60; K=6000; L=6000; M=6000; P=167 }
    136 = @{ D=44201; J=150; K=6500; L=6500; M=6500; P=181 }
    137 = @{ D=44208; J=180; K=6500; L=6500; M=6500; P=181 }
    138 = @{ D=44428; J=180; K=6500; L=6500; M=6500; P=181 }
    139 = @{ D=44463; J=180; K=6000; L=6000; M=6000; P=167 }
    140 = @{ D=44284; J=30; K=6000; L=6000; M=6000; P=167 }
    141 = @{ D=44210; J=60; K=6500; L=7000; M=6750; P=188 }
    142 = @{ D=44242; J=60; K=6000; L=6000; M=6000; P=167 }
    143 = @{ D=44518; J=60; K=5000; L=5000; M=5000; P=139 }
    144 = @{ D=44315; J=100; K=6000; L=6000; M=6000; P=167 }
    145 = @{ D=44370; J=32; K=6000; L=6000; M=6000; P=167 }
    146 = @{ D=44519; J=180; K=5000; L=5000; M=5000; P=139 }
    147 = @{ D=44462; J=90; K=6500; L=6500; M=6500; P=181 }
    148 = @{ D=44386; J=180; K=6000; L=6000; M=6000; P=167 }
    149 = @{ D=44168; J=60; K=6000; L=6000; M=6000; P=167 }
    150 = @{ D=44313; J=180; K=6000; L=6000; M=6000; P=167 }
    151 = @{ D=44435; J=500; K=6500; L=7000; M=6820; P=189 }
    152 = @{ D=44322; J=60; K=6000; L=6000; M=6000; P=167 }
    153 = @{ D=44231; J=60; K=6000; L=6000; M=6000; P=167 }
    154 = @{ D=44298; J=40; K=5500; L=6000; M=5750; P=160 }
    155 = @{ D=44357; J=60; K=6000; L=6000; M=6000; P=167 }
    156 = @{ D=44215; J=180; K=6000; L=6000; M=6000; P=167 }
    157 = @{ D=44204; J=150; K=6500; L=6500; M=6500; P=181 }
    158 = @{ D=44473; J=150; K=6000; L=6000; M=6000; P=167 }
    159 = @{ D=44200; J=70; K=7000; L=7000; M=7000; P=194 }
    160 = @{ D=44166; J=180; K=6000; L=6000; M=6000; P=167 }
    161 = @{ D=44259; J=60; K=6500; L=6500; M=6500; P=181 }
    162 = @{ D=44278; J=180; K=6000; L=6000; M=6000; P=167 }
    163 = @{ D=44522; J=60; K=5000; L=5000; M=5000; P=139 }
    164 = @{ D=44218; J=150; K=6000; L=6000; M=6000; P=167 }
    165 = @{ D=44396; J=60; K=6000; L=6000; M=6000; P=167 }
    166 = @{ D=44280; J=60; K=6000; L=6000; M=6000; P=167 }
    167 = @{ D=44427; J=70; K=7000; L=7000; M=7000; P=194 }
    168 = @{ D=44340; J=30; K=6000; L=6000; M=6000; P=167 }
    169 = @{ D=44497; J=100; K=5000; L=6000; M=5500; P=153 }
    170 = @{ D=44267; J=180; K=6000; L=6000; M=6000; P=167 }
    171 = @{ D=44418; J=180; K=7000; L=7000; M=7000; P=194 }
    172 = @{ D=44270; J=30; K=6000; L=6000; M=6000; P=167 }
    173 = @{ D=44525; J=70; K=6000; L=6000; M=6000; P=167 }
    174 = @{ D=44383; J=180; K=6000; L=6000; M=6000; P=167 }
    175 = @{ D=44286; J=27; K=6000; L=6000; M=6000; P=167 }
    176 = @{ D=44508; J=80; K=5500; L=5500; M=5500; P=153 }
    177 = @{ D=44307; J=30; K=6000; L=6000; M=6000; P=167 }
    178 = @{ D=44299; J=180; K=6000; L=6000; M=6000; P=167 }
    179 = @{ D=44316; J=180; K=6000; L=6000; M=6000; P=167 }
    180 = @{ D=44399; J=90; K=6000; L=6000; M=6000; P=167 }
    181 = @{ D=44229; J=180; K=6000; L=6000; M=6000; P=167 }
    182 = @{ D=44239; J=180; K=6000; L=6000; M=6000; P=167 }
    183 = @{ D=44469; J=90; K=6000; L=6000; M=6000; P=167 }
    184 = @{ D=44392; J=150; K=6000; L=6000; M=6000; P=167 }
    185 = @{ D=44160; J=22; K=6500; L=6500; M=6500; P=181 }
    186 = @{ D=44211; J=180; K=6000; L=6000; M=6000; P=167 }
    187 = @{ D=44273; J=30; K=6500; L=6500; M=6500; P=181 }
    188 = @{ D=44350; J=60; K=6000; L=6000; M=6000; P=167 }
    189 = @{ D=44358; J=140; K=6000; L=6000; M=6000; P=167 }
    190 = @{ D=44189; J=200; K=6000; L=7000; M=6500; P=181 }
    191 = @{ D=44335; J=30; K=6000; L=6000; M=6000; P=167 }
    192 = @{ D=44494; J=80; K=5500; L=5500; M=5500; P=153 }
    193 = @{ D=44342; J=28; K=6000; L=6000; M=6000; P=167 }
    194 = @{ D=44509; J=160; K=5000; L=5000; M=5000; P=139 }
    195 = @{ D=44491; J=180; K=5000; L=5500; M=5250; P=146 }
    196 = @{ D=44344; J=180; K=6000; L=6000; M=6000; P=167 }
    197 = @{ D=44474; J=180; K=5000; L=5000; M=5000; P=139 }
    198 = @{ D=44224; J=80; K=6000; L=6500; M=6250; P=174 }
    199 = @{ D=44447; J=36; K=6000; L=6000; M=6000; P=167 }
    200 = @{ D=44274; J=150; K=6500; L=6500; M=6500; P=181 }
    201 = @{ D=44490; J=60; K=5000; L=6000; M=5500; P=153 }
}

foreach ($r in $rowData.Keys) {
    $d = $rowData[$r]
    $ws.Cells.Item($r, 4).Value = $d.D
    $ws.Cells.Item($r, 10).Value = $d.J
    $ws.Cells.Item($r, 11).Value = $d.K
    $ws.Cells.Item($r, 12).Value = $d.L
    $ws.Cells.Item($r, 13).Value = $d.M
    $ws.Cells.Item($r, 16).Value = $d.P
}

$ws.Cells.Item(96, 15).Value = 'Provincia de Chacabuco'

# Add new row 202, cloned from former row 201 content
$ws.Cells.Item(202, 1).Value = 4
$ws.Cells.Item(202, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(202, 3).Value = 'Los Lagos'
$ws.Cells.Item(202, 4).Value = 44424
$ws.Cells.Item(202, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(202, 5).Value = 10
$ws.Cells.Item(202, 6).Value = 100112037
$ws.Cells.Item(202, 7).Value = 'Cebollín'
$ws.Cells.Item(202, 8).Value = 'Sin especificar'
$ws.Cells.Item(202, 9).Value = 'Primera'
$ws.Cells.Item(202, 10).Value = 60
$ws.Cells.Item(202, 11).Value = 7000
$ws.Cells.Item(202, 12).Value = 7000
$ws.Cells.Item(202, 13).Value = 7000
$ws.Cells.Item(202, 14).Value = '$/paquete 36 unidades'
$ws.Cells.Item(202, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(202, 16).Value = 194
$ws.Cells.Item(202, 17).Value = 36
$ws.Cells.Item(202, 18).Value = 'Hortaliza'
